# Generate Report for Handback
# Fills in the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" / "Handoff Reason" columns (E:H) for the two content rows on each
# language sheet, and updates the Status column to show the localized content
# has now been handed back and is in sync with en-US.

$wb = $excel.ActiveWorkbook

# Cornflowerblue (FF6495ED) expressed as the BGR-packed long that
# Range.Font.Color expects, matching the workbook's existing "HyperLink" look.
$hyperlinkColor = 15570276

function Set-HandbackCell {
    param($ws, $cellRef, $displayText, $url)

    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $displayText)
    $ws.Range($cellRef).Font.Underline = $true
    $ws.Range($cellRef).Font.Color = $hyperlinkColor
}

function Restore-LinkLook {
    param($ws, $cellRef)

    # Re-assert the existing "HyperLink" look (underline + cornflowerblue)
    # on a cell that already carries a hyperlink, so untouched link cells
    # keep their original appearance through the save round-trip.
    $ws.Range($cellRef).Font.Underline = $true
    $ws.Range($cellRef).Font.Color = $hyperlinkColor
}

function Set-HandbackRow {
    param($ws, $sheetName, $row, $handbackDateTime)

    $eCell = "E" + $row
    $fCell = "F" + $row
    $gCell = "G" + $row
    $hCell = "H" + $row

    $mdFile = $ws.Range("A" + $row).Value2
    $xlfFile = $ws.Range("C" + $row).Value2

    $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/b35c4b6463ac92267d65acaf504c120e6089c01f/e2e/" + $mdFile
    $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/" + $sheetName + "/" + $xlfFile

    # E = Latest Target File (same source doc that was handed off)
    Set-HandbackCell $ws $eCell $mdFile $mdUrl

    # F = Latest Handback File (the translated xlf that was returned)
    Set-HandbackCell $ws $fCell $xlfFile $xlfUrl

    # G = Latest Handback DateTime
    $ws.Range($gCell).Value = $handbackDateTime

    # H = Handoff Reason
    $ws.Range($hCell).Value = "Include"
}

$sheetNames = @("zh-cn", "de-de")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status column: handoff -> handed back, now in sync with en-US source.
    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    # Columns A/C (Source File Name, Latest Handoff File) already carry
    # hyperlinks from before this edit; keep their look intact.
    Restore-LinkLook $ws "A2"
    Restore-LinkLook $ws "C2"
    Restore-LinkLook $ws "A3"
    Restore-LinkLook $ws "C3"
    Restore-LinkLook $ws "A4"
}

# The Overview sheet mirrors each language's Status text in its own grid
# (column per language); keep it in step with the same wording change.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"
Restore-LinkLook $overview "A2"
Restore-LinkLook $overview "A3"
Restore-LinkLook $overview "A4"

$handbackDateTimeBySheet = @{
    "zh-cn" = "2016-03-11 05:12:13"
    "de-de" = "2016-03-11 05:12:39"
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $handbackDateTime = $handbackDateTimeBySheet[$sheetName]

    Set-HandbackRow $ws $sheetName 2 $handbackDateTime
    Set-HandbackRow $ws $sheetName 3 $handbackDateTime
}
